$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the three newly-documented journal entries (rows 73-75) ---
$ws.Range("C73").Value = "Fin de la documentation"
$ws.Range("D73").Value = 44287
$ws.Range("E73").Value = 150

$ws.Range("C74").Value = "Préparation de données factices finales"
$ws.Range("D74").Value = 44287
$ws.Range("E74").Value = 30

$ws.Range("C75").Value = "Mails de rendu"
$ws.Range("D75").Value = 44287
$ws.Range("E75").Value = 30

# Row 76 stays blank/unchanged.

# --- Remove the now-unused blank rows 77-90 (14 rows), shifting the ---
# --- trailing "spacer" row (old 91) and the bottom border row (old 92) up ---
$ws.Rows("77:90").Delete()

# --- Restore the current selection / scrolled view to match the author's ---
# --- last on-screen position when the sheet was saved ---
$ws.Range("E76").Select()
